$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair stats (Natmi following Dr Hou advice)
# Columns changed per row: E, G, H, I, J, K, M, N, O, P, Q, R, S, T
# (F and L remain unchanged)

$updates = @{
    2  = @{ E=3; G=6.072131; H=18.216393; I=0.003943999267036455; J=0.003943999267036454; K=3; M=112.513392; N=337.540176; O=0.3275312977368564; P=0.3275312977368564; Q=683.1960554783519; R=6148.764499305167; S=0.00129178319820566; T=0.00129178319820566 }
    3  = @{ E=3; G=6.072131; H=18.216393; I=0.003943999267036455; J=0.003943999267036454; K=3; M=106.314466; N=318.943398; O=0.3094859589441663; P=0.3094859589441664; Q=645.5553647470459; R=5809.998282723414; S=0.001220612395233866; T=0.001220612395233866 }
    4  = @{ E=3; G=6.072131; H=18.216393; I=0.003943999267036455; J=0.003943999267036454; K=3; M=124.6916553333333; N=374.074966; O=0.3629827433189773; P=0.3629827433189773; Q=757.1440657908487; R=6814.296592117638; S=0.001431603673596928; T=0.001431603673596928 }
    5  = @{ E=3; G=1480.851806666667; H=4442.55542; I=0.9618498744646554; J=0.9618498744646552; K=3; M=112.513392; N=337.540176; O=0.3275312977368564; P=0.3275312977368564; Q=166615.6598173949; R=1499540.938356554; S=0.3150359376114409; T=0.3150359376114409 }
    6  = @{ E=3; G=1480.851806666667; H=4442.55542; I=0.9618498744646554; J=0.9618498744646552; K=3; M=106.314466; N=318.943398; O=0.3094859589441663; P=0.3094859589441664; Q=157435.9690509019; R=1416923.721458117; S=0.2976790307590199; T=0.2976790307590199 }
    7  = @{ E=3; G=1480.851806666667; H=4442.55542; I=0.9618498744646554; J=0.9618498744646552; K=3; M=124.6916553333333; N=374.074966; O=0.3629827433189773; P=0.3629827433189773; Q=184649.863076624; R=1661848.767689616; S=0.3491349060941946; T=0.3491349060941945 }
    8  = @{ E=3; G=52.663316; H=157.989948; I=0.03420612626830831; J=0.0342061262683083; K=3; M=112.513392; N=337.540176; O=0.3275312977368564; P=0.3275312977368564; Q=5925.328317127872; R=53327.95485415084; S=0.01120357692720979; T=0.01120357692720979 }
    9  = @{ E=3; G=52.663316; H=157.989948; I=0.03420612626830831; J=0.0342061262683083; K=3; M=106.314466; N=318.943398; O=0.3094859589441663; P=0.3094859589441664; Q=5598.872318329256; R=50389.8508649633; S=0.01058631578991263; T=0.01058631578991263 }
    10 = @{ E=3; G=52.663316; H=157.989948; I=0.03420612626830831; J=0.0342061262683083; K=3; M=124.6916553333333; N=374.074966; O=0.3629827433189773; P=0.3629827433189773; Q=6566.67604738242; R=59100.08442644177; S=0.01241623355118588; T=0.01241623355118588 }
}

foreach ($r in $updates.Keys) {
    $rowVals = $updates[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
